$d = $word.ActiveDocument

# The opening sentence currently reads:
#   "Pensamento computacional: Processo de pensamento envolvido na expressao..."
# It needs to become:
#   "Pensamento computacional e o processo de pensamento envolvido na expressao..."
# split (as it was typed/edited interactively in Word) into the runs:
#   "Pensamento computacional e" | " o" | [[_GoBack bookmark]] | " p" | "rocesso de pensamento..."
#
# Paragraph 1 starts at document offset 0. The substring ": P" (colon, space,
# capital P) sits at offset 24-27 and must become " e o p" (accented e).

$r1 = $d.Range(24, 27)
$r1.Text = " é o p"

# At this point paragraph 1 is a single run:
#   "Pensamento computacional é o processo de pensamento envolvido..."
# Offsets (relative to the whole document, paragraph 1 still starts at 0):
#   26 -> right after "...computacional é"
#   28 -> right after "...é o"      (this is where the _GoBack bookmark goes)
#   30 -> right after "...é o p"    (start of "rocesso...")

# Split off run4 ("rocesso...") from run3 (" p") using a throwaway bookmark:
# adding (and then deleting) a bookmark at a position forces the run to be
# split there, and the split persists even after the bookmark is removed.
$d.Bookmarks.Add("ZZTempSplitA", $d.Range(30, 30))
$d.Bookmarks.Item("ZZTempSplitA").Delete()

# Place the real _GoBack bookmark between " o" and " p". Word only allows a
# single bookmark called "_GoBack" at a time, so re-adding it here also
# removes it from its old location (end of the "Design de algoritmos"
# paragraph).
$d.Bookmarks.Add("_GoBack", $d.Range(28, 28))

# Split off run2 (" o") from run1 ("Pensamento computacional é") the same way.
$d.Bookmarks.Add("ZZTempSplitB", $d.Range(26, 26))
$d.Bookmarks.Item("ZZTempSplitB").Delete()
